$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.Formula = "'248.67"
$c.Style = $s
$c = $ws.Range("G2")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D3")
$s = $c.Style
$c.Formula = "'21.69"
$c.Style = $s
$c = $ws.Range("G3")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D4")
$s = $c.Style
$c.Formula = "'5.419"
$c.Style = $s
$c = $ws.Range("G4")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D5")
$s = $c.Style
$c.Formula = "'0.05692"
$c.Style = $s
$c = $ws.Range("G5")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G6")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D7")
$s = $c.Style
$c.Formula = "'0.8067"
$c.Style = $s
$c = $ws.Range("G7")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D8")
$s = $c.Style
$c.Formula = "'1.033"
$c.Style = $s
$c = $ws.Range("G8")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$c = $ws.Range("D9")
$s = $c.Style
$c.Formula = "'0.01168"
$c.Style = $s
$ws.Range("E9").Value = "8OneONE"
$c = $ws.Range("G9")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Range("D10")
$s = $c.Style
$c.Formula = "'0.1462"
$c.Style = $s
$ws.Range("E10").Value = "9WazirXWRX"
$c = $ws.Range("G10")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Range("D11")
$s = $c.Style
$c.Formula = "'0.07717"
$c.Style = $s
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$c = $ws.Range("G11")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c = $ws.Range("D12")
$s = $c.Style
$c.Formula = "'0.03195"
$c.Style = $s
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$c = $ws.Range("G12")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Range("D13")
$s = $c.Style
$c.Formula = "'0.03073"
$c.Style = $s
$ws.Range("E13").Value = "12BitrueCoinBTR"
$c = $ws.Range("G13")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Range("D14")
$s = $c.Style
$c.Formula = "'0.09275"
$c.Style = $s
$ws.Range("E14").Value = "13BitMartTokenBMX"
$c = $ws.Range("G14")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$c = $ws.Range("D15")
$s = $c.Style
$c.Formula = "'3.558"
$c.Style = $s
$ws.Range("E15").Value = "14MCDexMCB"
$c = $ws.Range("G15")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Range("D16")
$s = $c.Style
$c.Formula = "'0.001646"
$c.Style = $s
$ws.Range("E16").Value = "15BitForexTokenBF"
$c = $ws.Range("G16")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$c = $ws.Range("D17")
$s = $c.Style
$c.Formula = "'0.04709"
$c.Style = $s
$ws.Range("E17").Value = "16CoinExTokenCET"
$c = $ws.Range("G17")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D18")
$s = $c.Style
$c.Formula = "'0.006361"
$c.Style = $s
$c = $ws.Range("G18")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D19")
$s = $c.Style
$c.Formula = "'0.005049"
$c.Style = $s
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"
$c = $ws.Range("G19")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D20")
$s = $c.Style
$c.Formula = "'0.001043"
$c.Style = $s
$c = $ws.Range("G20")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D21")
$s = $c.Style
$c.Formula = "'0.0001501"
$c.Style = $s
$c = $ws.Range("G21")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D22")
$s = $c.Style
$c.Formula = "'0.0003201"
$c.Style = $s
$c = $ws.Range("G22")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D23")
$s = $c.Style
$c.Formula = "'3.773"
$c.Style = $s
$c = $ws.Range("G23")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D24")
$s = $c.Style
$c.Formula = "'6.429"
$c.Style = $s
$c = $ws.Range("G24")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D25")
$s = $c.Style
$c.Formula = "'2.165"
$c.Style = $s
$c = $ws.Range("G25")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G26")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D27")
$s = $c.Style
$c.Formula = "'0.1318"
$c.Style = $s
$c = $ws.Range("G27")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G28")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G29")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G30")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G31")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G32")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G33")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G34")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G35")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G36")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G37")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G38")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G39")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D40")
$s = $c.Style
$c.Formula = "'0.04125"
$c.Style = $s
$c = $ws.Range("G40")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D41")
$s = $c.Style
$c.Formula = "'0.006947"
$c.Style = $s
$ws.Range("E41").Value = "40KickTokenKICK"
$c = $ws.Range("G41")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$c = $ws.Range("D42")
$s = $c.Style
$c.Formula = "'0.003501"
$c.Style = $s
$ws.Range("E42").Value = "41CEJICEJI"
$c = $ws.Range("G42")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$c = $ws.Range("D43")
$s = $c.Style
$c.Formula = "'0.1044"
$c.Style = $s
$ws.Range("E43").Value = "42BKEXTokenBKK"
$c = $ws.Range("G43")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D44")
$s = $c.Style
$c.Formula = "'0.007952"
$c.Style = $s
$c = $ws.Range("G44")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D45")
$s = $c.Style
$c.Formula = "'0.00005905"
$c.Style = $s
$c = $ws.Range("G45")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G46")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D47")
$s = $c.Style
$c.Formula = "'0.0005503"
$c.Style = $s
$c = $ws.Range("G47")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D48")
$s = $c.Style
$c.Formula = "'0.6827"
$c.Style = $s
$c = $ws.Range("G48")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D49")
$s = $c.Style
$c.Formula = "'0.009068"
$c.Style = $s
$c = $ws.Range("G49")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("D50")
$s = $c.Style
$c.Formula = "'0.00002101"
$c.Style = $s
$c = $ws.Range("G50")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
$c = $ws.Range("G51")
$s = $c.Style
$c.Formula = "'13"
$c.Style = $s
